$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---- Row 6: question #5 - "角色管理" (Role management) ----
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "角色管理"
$ws.Range("C6").Value = "角色对应权限列表里的移除权限的方法没有提供。"
$ws.Range("D6").Value = "耿晓红"
# "2015.11.12" looks like a date, so Excel would normally auto-convert a
# typed value into a date serial number. Enter it as a text formula first,
# then bake it down to a plain (shared-string) value so it stays textual.
$ws.Range("E6").Formula = "=""2015.11.12"""
$ws.Range("E6").Copy()
$ws.Range("E6").PasteSpecial(-4163)
$ws.Rows.Item(6).RowHeight = 27

# ---- Row 7: question #6 - "权限管理" (Permission management) ----
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "权限管理"
$ws.Range("C7").Value = "在给的需求文档里没有关于权限管理的修改的接口方法，请提供。"
$ws.Range("D7").Value = "耿晓红"
$ws.Range("E7").Formula = "=""2015.11.12"""
$ws.Range("E7").Copy()
$ws.Range("E7").PasteSpecial(-4163)
$ws.Rows.Item(7).RowHeight = 27

# ---- Update the selected / visible range in the sheet view ----
$ws.Range("D7:E7").Select()
